# [12.12.2023][add Accepted with warning flow]
#
# TestData.xlsx - add the "Accepted with warning" CSID test-case rows
# (Standard / Simplified / Both) between the existing "Accepted" CSID
# rows and the "Accepted Clearance/Reporting" rows, rename the existing
# CSID test cases to the "...Accepted..." variants, and add a 4th
# Clearance sample row (BR-KSA-20_Standard_Note.xml).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1) Make room: copy the *formatting* of the existing template rows into
#    the new rows first (values are copied later), using Formats-only
#    paste so nothing here disturbs the still-present old values we
#    need as templates for further copies.
# ---------------------------------------------------------------------

# 7-column "Accepted" header/value row pair -> reused for the new
# Warning header/value row pairs (rows 8/9 and 10/11).
$ws.Range("A2:G3").Copy()
$ws.Range("A8:G9").PasteSpecial(-4122)
$ws.Range("A10:G11").PasteSpecial(-4122)

# 10-column "Both" header/value row pair -> reused for the new "Both
# Warning" header/value row pair (rows 12/13).
$ws.Range("A6:J7").Copy()
$ws.Range("A12:J13").PasteSpecial(-4122)

# 4-column Clearance header/value rows -> reused for the extra
# Clearance sample row and for the whole Reporting block, which both
# move down from their old positions (old rows 8-15 -> new rows
# 14-22).
$ws.Range("A8:D11").Copy()
$ws.Range("A14:D17").PasteSpecial(-4122)
$ws.Range("A12:D15").Copy()
$ws.Range("A18:D21").PasteSpecial(-4122)
$ws.Range("A12:D12").Copy()
$ws.Range("A22").PasteSpecial(-4122)

$excel.CutCopyMode = 0

# ---------------------------------------------------------------------
# 2) Write the cell values. Formatting (incl. the Text number format
#    used for the vatNumber column) is already in place, so the
#    vatNumber digit-strings are stored as text, matching the source.
# ---------------------------------------------------------------------

# Row 2-3 : testCSIDStandardAcceptedInvoice (renamed from ...StandardInvoice)
$ws.Range("A2").Value = "testCSIDStandardAcceptedInvoice"
$ws.Range("C2").Value = "invoiceFileName"
$ws.Range("D2").Value = "invoiceCreditFileName"
$ws.Range("E2").Value = "invoiceDebitFileName"
$ws.Range("F2").Value = "csrFileName"
$ws.Range("G2").Value = "vatNumber"

$ws.Range("A3").Value = "testCSIDStandardAcceptedInvoice"
$ws.Range("B3").Value = $true
$ws.Range("C3").Value = "Standard_Invoice.xml"
$ws.Range("D3").Value = "Standard_Credit_Note.xml"
$ws.Range("E3").Value = "Standard_Debit_Note.xml"
$ws.Range("F3").Value = "csr-config-example-EN.properties"
$ws.Range("G3").Value = "310094010300003"

# Row 4-5 : testCSIDSimplifiedAcceptedInvoice (renamed from ...SimplifiedInvoice)
$ws.Range("A4").Value = "testCSIDSimplifiedAcceptedInvoice"
$ws.Range("C4").Value = "invoiceFileName"
$ws.Range("D4").Value = "invoiceCreditFileName"
$ws.Range("E4").Value = "invoiceDebitFileName"
$ws.Range("F4").Value = "csrFileName"
$ws.Range("G4").Value = "vatNumber"

$ws.Range("A5").Value = "testCSIDSimplifiedAcceptedInvoice"
$ws.Range("B5").Value = $true
$ws.Range("C5").Value = "Simplified_Invoice.xml"
$ws.Range("D5").Value = "Simplified_Credit_Note.xml"
$ws.Range("E5").Value = "Simplified_Debit_Note.xml"
$ws.Range("F5").Value = "csr-config-example-EN-simplified.properties"
$ws.Range("G5").Value = "310094010300003"

# Row 6-7 : testCSIDBothAcceptedInvoice (renamed from ...BothInvoice)
$ws.Range("A6").Value = "testCSIDBothAcceptedInvoice"
$ws.Range("C6").Value = "invoiceFileName"
$ws.Range("D6").Value = "invoiceCreditFileName"
$ws.Range("E6").Value = "invoiceDebitFileName"
$ws.Range("F6").Value = "standardInvoiceFileName"
$ws.Range("G6").Value = "standardInvoiceCreditFileName"
$ws.Range("H6").Value = "standardInvoiceDebitFileName"
$ws.Range("I6").Value = "csrFileName"
$ws.Range("J6").Value = "vatNumber"

$ws.Range("A7").Value = "testCSIDBothAcceptedInvoice"
$ws.Range("B7").Value = $true
$ws.Range("C7").Value = "Simplified_Invoice.xml"
$ws.Range("D7").Value = "Simplified_Credit_Note.xml"
$ws.Range("E7").Value = "Simplified_Debit_Note.xml"
$ws.Range("F7").Value = "Standard_Invoice.xml"
$ws.Range("G7").Value = "Standard_Credit_Note.xml"
$ws.Range("H7").Value = "Standard_Debit_Note.xml"
$ws.Range("I7").Value = "csr-config-example-EN-both.properties"
$ws.Range("J7").Value = "310094010300003"

# Row 8-9 : testCSIDStandardWarningInvoice (new)
$ws.Range("A8").Value = "testCSIDStandardWarningInvoice"
$ws.Range("C8").Value = "invoiceFileName"
$ws.Range("D8").Value = "invoiceCreditFileName"
$ws.Range("E8").Value = "invoiceDebitFileName"
$ws.Range("F8").Value = "csrFileName"
$ws.Range("G8").Value = "vatNumber"

$ws.Range("A9").Value = "testCSIDStandardWarningInvoice"
$ws.Range("B9").Value = $true
$ws.Range("C9").Value = "BR-KSA-35_BR-KSA-15.xml"
$ws.Range("D9").Value = "BR-KSA-36.xml"
$ws.Range("E9").Value = "BR-KSA-36.xml"
$ws.Range("F9").Value = "csr-config-example-EN.properties"
$ws.Range("G9").Value = "310094010300003"

# Row 10-11 : testCSIDSimplifiedWarningInvoice (new)
$ws.Range("A10").Value = "testCSIDSimplifiedWarningInvoice"
$ws.Range("C10").Value = "invoiceFileName"
$ws.Range("D10").Value = "invoiceCreditFileName"
$ws.Range("E10").Value = "invoiceDebitFileName"
$ws.Range("F10").Value = "csrFileName"
$ws.Range("G10").Value = "vatNumber"

$ws.Range("A11").Value = "testCSIDSimplifiedWarningInvoice"
$ws.Range("B11").Value = $true
$ws.Range("C11").Value = "BR-KSA-35.xml"
$ws.Range("D11").Value = "BR-KSA-83_BR-KSA-F-06-C16.xml"
$ws.Range("E11").Value = "BR-KSA-56.xml"
$ws.Range("F11").Value = "csr-config-example-EN-simplified.properties"
$ws.Range("G11").Value = "310094010300003"

# Row 12-13 : testCSIDBothWarningInvoice (new)
$ws.Range("A12").Value = "testCSIDBothWarningInvoice"
$ws.Range("C12").Value = "invoiceFileName"
$ws.Range("D12").Value = "invoiceCreditFileName"
$ws.Range("E12").Value = "invoiceDebitFileName"
$ws.Range("F12").Value = "standardInvoiceFileName"
$ws.Range("G12").Value = "standardInvoiceCreditFileName"
$ws.Range("H12").Value = "standardInvoiceDebitFileName"
$ws.Range("I12").Value = "csrFileName"
$ws.Range("J12").Value = "vatNumber"

$ws.Range("A13").Value = "testCSIDBothWarningInvoice"
$ws.Range("B13").Value = $true
$ws.Range("C13").Value = "BR-KSA-35.xml"
$ws.Range("D13").Value = "BR-KSA-83_BR-KSA-F-06-C16.xml"
$ws.Range("E13").Value = "BR-KSA-56.xml"
$ws.Range("F13").Value = "BR-KSA-35_BR-KSA-15.xml"
$ws.Range("G13").Value = "BR-KSA-36.xml"
$ws.Range("H13").Value = "BR-KSA-36.xml"
$ws.Range("I13").Value = "csr-config-example-EN-both.properties"
$ws.Range("J13").Value = "310094010300003"

# Row 14-18 : testAcceptedClearanceInvoice (moved down, +1 sample row)
$ws.Range("A14").Value = "testAcceptedClearanceInvoice"
$ws.Range("C14").Value = "invoiceType"
$ws.Range("D14").Value = "invoiceFileName"

$ws.Range("A15").Value = "testAcceptedClearanceInvoice"
$ws.Range("B15").Value = $true
$ws.Range("C15").Value = "STANDARDNOTE"
$ws.Range("D15").Value = "Standard_Invoice.xml"

$ws.Range("A16").Value = "testAcceptedClearanceInvoice"
$ws.Range("B16").Value = $true
$ws.Range("C16").Value = "STANDARDCREDIT"
$ws.Range("D16").Value = "Standard_Credit_Note.xml"

$ws.Range("A17").Value = "testAcceptedClearanceInvoice"
$ws.Range("B17").Value = $true
$ws.Range("C17").Value = "STANDARDDEBIT"
$ws.Range("D17").Value = "Standard_Debit_Note.xml"

$ws.Range("A18").Value = "testAcceptedClearanceInvoice"
$ws.Range("B18").Value = $true
$ws.Range("C18").Value = "STANDARDNOTE"
$ws.Range("D18").Value = "BR-KSA-20_Standard_Note.xml"

# Row 19-22 : testAcceptedReportingInvoice (moved down)
$ws.Range("A19").Value = "testAcceptedReportingInvoice"
$ws.Range("C19").Value = "invoiceType"
$ws.Range("D19").Value = "invoiceFileName"

$ws.Range("A20").Value = "testAcceptedReportingInvoice"
$ws.Range("B20").Value = $true
$ws.Range("C20").Value = "SIMPLIFIEDNOTE"
$ws.Range("D20").Value = "Simplified_Invoice.xml"

$ws.Range("A21").Value = "testAcceptedReportingInvoice"
$ws.Range("B21").Value = $true
$ws.Range("C21").Value = "SIMPLIFIEDCREDIT"
$ws.Range("D21").Value = "Simplified_Credit_Note.xml"

$ws.Range("A22").Value = "testAcceptedReportingInvoice"
$ws.Range("B22").Value = $true
$ws.Range("C22").Value = "SIMPLIFIEDDEBIT"
$ws.Range("D22").Value = "Simplified_Debit_Note.xml"

# ---------------------------------------------------------------------
# 3) View state: the active selection moved from B13 to A13, and the
#    window scrolled down a little (top-left visible row ~A7).
# ---------------------------------------------------------------------
[void]$ws.Range("A13").Select()
$excel.ActiveWindow.ScrollRow = 7
$excel.ActiveWindow.ScrollColumn = 1
